# The source workbook for year 2019 was repurposed for year 2005: the
# lone data row (A2) now holds the 2005 reporting year, and the sheet's
# cursor/selection moves to that new cell - mirroring what Excel records
# when a user types a value into A2 and leaves the selection there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new data row.
$ws.Range("A2").Value = 2005

# Move the active selection to the newly entered cell (updates
# dimension/selection in the saved sheetView, same as Excel would).
$ws.Range("A2").Select()
